$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "average_color" column (L) with a header and four data rows.
# The shared-string table must pick up "average_color", "160,104,73",
# "220,55,70", "105,76,77", "243,131,27" in that exact order, so the cells
# are written in that order (L1, L2, L4, L3, L5).
$ws.Range("L1").Value = "average_color"
$ws.Range("L2").Value = "160,104,73"
$ws.Range("L4").Value = "220,55,70"
$ws.Range("L3").Value = "105,76,77"
$ws.Range("L5").Value = "243,131,27"

# Match the header cell L1 with the same direct formatting already used by
# the J1 header cell (bold font + side borders + centered alignment).
$ws.Range("J1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Give column L (12) a custom width close to the target 14.28515625.
$ws.Columns(12).ColumnWidth = 13.5

# Update the active selection to L6, matching the saved selection in the file.
$ws.Range("L6").Select() | Out-Null
